{"js": "// Update the answer table: each cell below is addressed by its\n// (row, col) position in the first table and its text content is\n// replaced with the new value (formatting/run-properties untouched).\nconst updates = [\n  { row: 0, col: 0, text: '80\u00f79=8, 8' },\n  { row: 0, col: 1, text: '42\u00f78=5, 2' },\n  { row: 0, col: 2, text: '89\u00f72=44, 1' },\n  { row: 0, col: 3, text: '27\u00f74=6, 3' },\n  { row: 0, col: 4, text: '20\u00f74=5, 0' },\n  { row: 4, col: 0, text: '78\u00f73=26, 0' },\n  { row: 4, col: 1, text: '35\u00f75=7, 0' },\n  { row: 4, col: 2, text: '20\u00f77=2, 6' },\n  { row: 4, col: 4, text: '32\u00f77=4, 4' },\n  { row: 8, col: 0, text: '68\u00f72=34, 0' },\n  { row: 8, col: 1, text: '48\u00f77=6, 6' },\n  { row: 8, col: 2, text: '20\u00f73=6, 2' },\n  { row: 8, col: 3, text: '32\u00f78=4, 0' },\n  { row: 8, col: 4, text: '69\u00f74=17, 1' },\n  { row: 12, col: 0, text: '50\u00f73=16, 2' },\n  { row: 12, col: 1, text: '35\u00f74=8, 3' },\n  { row: 12, col: 2, text: '43\u00f72=21, 1' },\n  { row: 12, col: 3, text: '48\u00f72=24, 0' },\n  { row: 12, col: 4, text: '84\u00f76=14, 0' },\n  { row: 16, col: 0, text: '73\u00f78=9, 1' },\n  { row: 16, col: 1, text: '21\u00f72=10, 1' },\n  { row: 16, col: 2, text: '94\u00f77=13, 3' },\n  { row: 16, col: 3, text: '65\u00f79=7, 2' },\n  { row: 16, col: 4, text: '20\u00f73=6, 2' },\n];\n\nconst tables = context.document.body.tables;\ntables.load('items');\nawait context.sync();\n\nconst table = tables.items[0];\n\nfor (const u of updates) {\n  const cell = table.getCell(u.row, u.col);\n  cell.value = u.text;\n}\n\nawait context.sync();\n", "ps1": "# Update the answer table: each entry below is addressed by its\n# 1-based (Row, Col) position in the document's first table, and its\n# text is replaced with the new value. Assigning to Range.Text leaves\n# the end-of-cell marker and the run's formatting (font/size) intact.\n$d = $word.ActiveDocument\n$t = $d.Tables.Item(1)\n\n$updates = @(\n    @{ Row = 1; Col = 1; Text = '80\u00f79=8, 8' },\n    @{ Row = 1; Col = 2; Text = '42\u00f78=5, 2' },\n    @{ Row = 1; Col = 3; Text = '89\u00f72=44, 1' },\n    @{ Row = 1; Col = 4; Text = '27\u00f74=6, 3' },\n    @{ Row = 1; Col = 5; Text = '20\u00f74=5, 0' },\n    @{ Row = 5; Col = 1; Text = '78\u00f73=26, 0' },\n    @{ Row = 5; Col = 2; Text = '35\u00f75=7, 0' },\n    @{ Row = 5; Col = 3; Text = '20\u00f77=2, 6' },\n    @{ Row = 5; Col = 5; Text = '32\u00f77=4, 4' },\n    @{ Row = 9; Col = 1; Text = '68\u00f72=34, 0' },\n    @{ Row = 9; Col = 2; Text = '48\u00f77=6, 6' },\n    @{ Row = 9; Col = 3; Text = '20\u00f73=6, 2' },\n    @{ Row = 9; Col = 4; Text = '32\u00f78=4, 0' },\n    @{ Row = 9; Col = 5; Text = '69\u00f74=17, 1' },\n    @{ Row = 13; Col = 1; Text = '50\u00f73=16, 2' },\n    @{ Row = 13; Col = 2; Text = '35\u00f74=8, 3' },\n    @{ Row = 13; Col = 3; Text = '43\u00f72=21, 1' },\n    @{ Row = 13; Col = 4; Text = '48\u00f72=24, 0' },\n    @{ Row = 13; Col = 5; Text = '84\u00f76=14, 0' },\n    @{ Row = 17; Col = 1; Text = '73\u00f78=9, 1' },\n    @{ Row = 17; Col = 2; Text = '21\u00f72=10, 1' },\n    @{ Row = 17; Col = 3; Text = '94\u00f77=13, 3' },\n    @{ Row = 17; Col = 4; Text = '65\u00f79=7, 2' },\n    @{ Row = 17; Col = 5; Text = '20\u00f73=6, 2' }\n)\n\nforeach ($u in $updates) {\n    $cell = $t.Cell($u.Row, $u.Col)\n    $cell.Range.Text = $u.Text\n}\n"}
